$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New donation rows (rows 4 through 12) captured from additional "Send Email OTP"
# submissions. Columns: A Receipt Number, B Row Number, C Timestamp, D First Name,
# E Last Name, F Donation Amount, G Email, H Phone, I Street Address, J City,
# K State, L Zip.
# Columns F, H and L look numeric but must be stored as text (matching the
# existing rows 2-3), so they are entered with a leading apostrophe to force
# Excel to treat them as text instead of auto-converting to numbers.
$data = @(
    @("REC-1741876652886-252", 3,  "2025-03-13T14:37:32.907Z", "Deepak", "Adhikari", "'3445", "dadhikari856@gmail.com", "'3477712375", "11 alpine ln", "Hicksville", "NY", "'11801"),
    @("REC-1741876747020-229", 4,  "2025-03-13T14:39:07.031Z", "Deepak", "Adhikari", "'3445", "dadhikari856@gmail.com", "'3477712375", "11 alpine ln", "Hicksville", "NY", "'11801"),
    @("REC-1741876748126-582", 5,  "2025-03-13T14:39:08.134Z", "Deepak", "Adhikari", "'3445", "dadhikari856@gmail.com", "'3477712375", "11 alpine ln", "Hicksville", "NY", "'11801"),
    @("REC-1741876748336-135", 6,  "2025-03-13T14:39:08.338Z", "Deepak", "Adhikari", "'3445", "dadhikari856@gmail.com", "'3477712375", "11 alpine ln", "Hicksville", "NY", "'11801"),
    @("REC-1741876873473-419", 7,  "2025-03-13T14:41:13.475Z", "Deepak", "Adhikari", "'3445", "dadhikari856@gmail.com", "'3477712375", "11 alpine ln", "Hicksville", "NY", "'11801"),
    @("REC-1741876962971-397", 8,  "2025-03-13T14:42:42.973Z", "Deepak", "Adhikari", "'3445", "dadhikari856@gmail.com", "'3477712375", "11 alpine ln", "Hicksville", "NY", "'11801"),
    @("REC-1741876992736-933", 9,  "2025-03-13T14:43:12.740Z", "Deepak", "Adhikari", "'3445", "dadhikari856@gmail.com", "'3477712375", "11 alpine ln", "Hicksville", "NY", "'11801"),
    @("REC-1741877009123-108", 10, "2025-03-13T14:43:29.132Z", "Deepak", "Adhikari", "'3445", "dadhikari856@gmail.com", "'3477712375", "11 alpine ln", "Hicksville", "NY", "'11801"),
    @("REC-1741877081275-061", 11, "2025-03-13T14:44:41.278Z", "Deepak", "Adhikari", "'3445", "dadhikari856@gmail.com", "'3477712375", "11 alpine ln", "Hicksville", "NY", "'11801")
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowData[$col - 1]
    }
}
